# Actualización automática 2025-10-07 16:30:08
$wb = $excel.ActiveWorkbook

# ---- Sheet "VENTAS POR GRUPO" ----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D9").Value = 915.84
$ws1.Range("E9").Value = 811.24
$ws1.Range("I9").Value = 226.8
$ws1.Range("E52").Value = 443.35
$ws1.Range("D56").Value = "2 de 54"
$ws1.Range("E56").Value = "2 de 54"
$ws1.Range("I56").Value = "4 de 54"

# ---- Sheet "VENTA MENSUAL" ----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F9").Value = 1953.88
$ws2.Range("F53").Value = 443.35
$ws2.Range("F54").Value = 443.35
$ws2.Range("F60").Value = 7759.8

# ---- Sheet "CUMPLIMIENTO MENSUAL" ----
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 1373.76
$ws3.Range("E3").Value = 16295.3870988183
$ws3.Range("F3").Value = 0.0777490838871264

$ws3.Range("D4").Value = 1254.59
$ws3.Range("E4").Value = -211.36711473472
$ws3.Range("F4").Value = 1.202609737305534

$ws3.Range("D7").Value = 420.3
$ws3.Range("E7").Value = 466.411016287574
$ws3.Range("F7").Value = 0.4739988477414949

$ws3.Range("D14").Value = 7275.23
$ws3.Range("E14").Value = 91741.27661190613
$ws3.Range("F14").Value = 0.07347492098984229
